$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.26
$ws.Range("B3").Value = 0.13
$ws.Range("B4").Value = -0.15
$ws.Range("B5").Value = -0.39
$ws.Range("B6").Value = 0.04
$ws.Range("B7").Value = -0.37
$ws.Range("B8").Value = -0.01
$ws.Range("B9").Value = 0.01
$ws.Range("B10").Value = 0.1
$ws.Range("B11").Value = 0.09
$ws.Range("B12").Value = 0.11
$ws.Range("B13").Value = -0.13
